$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style used for the other
# header cells (bold font, thin border all around, centered/top aligned).
# Copy the formatting from an existing header cell rather than re-applying
# individual properties, so the same style is reused instead of a new one
# being generated.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I (I0) and J (IF) columns, rows 2-25
$data = @{
    2  = @(6, 6)
    3  = @(7, 8)
    4  = @(7, 7)
    5  = @(8, 8)
    6  = @(7, 7)
    7  = @(8, 9)
    8  = @(8, 8)
    9  = @(7, 7)
    10 = @(6, 7)
    11 = @(8, 8)
    12 = @(7, 8)
    13 = @(7, 8)
    14 = @(8, 9)
    15 = @(6, 6)
    16 = @(6, 7)
    17 = @(9, 9)
    18 = @(8, 8)
    19 = @(8, 8)
    20 = @(9, 9)
    21 = @(9, 9)
    22 = @(8, 8)
    23 = @(6, 6)
    24 = @(4, 4)
    25 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
